{"js": "// UC5 Ret Enhed.docx fixes:\n//  1) \"[Udvidelse 1: Bruger v\u00e6lger intet felt.]\" -> \"[Udvidelse 1: Bruger annullerer indtastningen.]\"\n//  2) Move the \"_GoBack\" bookmark from the end of \"Forts\u00e6t i punkt 5.\" to the\n//     middle of the word \"punkt\" (between \"pun\" and \"kt\").\n\nconst body = context.document.body;\n\n// --- 1) Fix the extension-1 sentence under \"Udvidelser\" --------------------\nconst target = body.search(\"Bruger v\u00e6lger intet felt.\", { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\"Bruger annullerer indtastningen.\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2) Move the \"_GoBack\" bookmark into the middle of \"punkt\" -------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst prefix = body.search(\"Forts\u00e6t i pun\", { matchCase: true });\nprefix.load(\"text\");\nawait context.sync();\n\nif (prefix.items.length > 0) {\n  const splitPoint = prefix.items[0].getRange(\"End\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# UC5 Ret Enhed.docx fixes:\n#  1) \"[Udvidelse 1: Bruger v\u00e6lger intet felt.]\" -> \"[Udvidelse 1: Bruger annullerer indtastningen.]\"\n#  2) Move the \"_GoBack\" bookmark from the end of \"Forts\u00e6t i punkt 5.\" to the\n#     middle of the word \"punkt\" (between \"pun\" and \"kt\").\n\n$d = $word.ActiveDocument\n\n# --- 1) Fix the extension-1 sentence under \"Udvidelser\" -------------------\n$needle = \"Bruger v\u00e6lger intet felt.\"\n$text = $d.Content.Text\n$idx = $text.IndexOf($needle)\nif ($idx -ge 0) {\n    $rng = $d.Range($idx, $idx + $needle.Length)\n    $rng.Text = \"Bruger annullerer indtastningen.\"\n}\n\n# --- 2) Move the \"_GoBack\" bookmark into the middle of \"punkt\" ------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$text = $d.Content.Text\n$idx = $text.IndexOf(\"Forts\u00e6t i punkt 5\")\nif ($idx -ge 0) {\n    $splitPos = $idx + (\"Forts\u00e6t i pun\").Length\n    $target = $d.Range($splitPos, $splitPos)\n    $d.Bookmarks.Add(\"_GoBack\", $target)\n}\n"}
